$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2 = @{ "B" = 1.119349272945897; "C" = 0.4034454972101571; "D" = 0.5750552436731056; "E" = 0.2058480855811098; "G" = 1.827033437555315; "H" = 1.54303804516428; "I" = 1.186425040851958; "J" = 0.08799914823944732; "M" = 0.5851876514217054 }
    3 = @{ "B" = 1.033197247967792; "C" = 0.3728568879844829; "D" = 0.5717112469103398; "E" = 0.2055798808092426; "G" = 1.804951866953104; "H" = 1.539796088774011; "I" = 1.18544203092933; "J" = 0.08837240945389979; "M" = 0.5606221680284804 }
    4 = @{ "B" = 0.9807837804196708; "C" = 0.3542735153050671; "D" = 0.5699349536108542; "E" = 0.2055064897186725; "G" = 1.792638362165235; "H" = 1.538650118502545; "I" = 1.185551612859108; "J" = 0.08863323954474822; "M" = 0.5458576528635533 }
    5 = @{ "B" = 0.959546690127894; "C" = 0.3467501998440241; "D" = 0.5692806657377787; "E" = 0.2054995163300326; "G" = 1.787931781773409; "H" = 1.538394868596384; "I" = 1.185774838267292; "J" = 0.0887474828028143; "M" = 0.539921046419785 }
    6 = @{ "B" = 0.956027647786641; "C" = 0.3455039447201216; "D" = 0.5691762213363063; "E" = 0.2054997429702681; "G" = 1.78716900795942; "H" = 1.538365254622818; "I" = 1.185822667714433; "J" = 0.08876693300618932; "M" = 0.5389401094938009 }
    7 = @{ "B" = 0.9804968762364012; "C" = 0.354171852920615; "D" = 0.5699258480762239; "E" = 0.2055063028409982; "G" = 1.792573629515914; "H" = 1.538645819640664; "I" = 1.185553901348356; "J" = 0.08863474807490057; "M" = 0.5457772657182005 }
    8 = @{ "B" = 1.089543553351064; "C" = 0.3928572427655581; "D" = 0.5738447097614738; "E" = 0.2057366484323389; "G" = 1.819160339106503; "H" = 1.541744529686895; "I" = 1.185937623667989; "J" = 0.08812127753895638; "M" = 0.576651230235953 }
    9 = @{ "B" = 1.307236679393668; "C" = 0.4703082339450475; "D" = 0.5837313570833942; "E" = 0.2069139140343417; "G" = 1.881254840693316; "H" = 1.554553836002299; "I" = 1.192384365667493; "J" = 0.0873657171523945; "M" = 0.6397336476156426 }
    10 = @{ "B" = 1.469559267126613; "C" = 0.5282101694770063; "D" = 0.5923452087204453; "E" = 0.2082232914172728; "G" = 1.933067419856428; "H" = 1.568115832336673; "I" = 1.200643766092185; "J" = 0.0869642487431932; "M" = 0.6876468639535034 }
    11 = @{ "B" = 1.543930685318628; "C" = 0.5547751472187201; "D" = 0.596558902110047; "E" = 0.2089159658392141; "G" = 1.958009479760875; "H" = 1.575197231769948; "I" = 1.205177483329209; "J" = 0.08681506857832488; "M" = 0.709788603342659 }
    12 = @{ "B" = 1.572169822101102; "C" = 0.5648674333567669; "D" = 0.5981970942262649; "E" = 0.2091922518931923; "G" = 1.967653760780166; "H" = 1.578010708386216; "I" = 1.207006829233663; "J" = 0.08676339573165492; "M" = 0.7182230942174925 }
    13 = @{ "B" = 1.566084623994698; "C" = 0.5626924188763383; "D" = 0.5978423857252153; "E" = 0.2091321262256258; "G" = 1.965567798986797; "H" = 1.577398897134316; "I" = 1.206607828342655; "J" = 0.08677430996624835; "M" = 0.7164043537543847 }
    14 = @{ "B" = 1.54625240650347; "C" = 0.5556047882297435; "D" = 0.5966928236109368; "E" = 0.2089384155853011; "G" = 1.958798915872194; "H" = 1.575426050226753; "I" = 1.20532572453611; "J" = 0.08681072080443997; "M" = 0.7104815129928852 }
    15 = @{ "B" = 1.534114539194547; "C" = 0.5512676792530442; "D" = 0.5959942291586628; "E" = 0.2088215844940393; "G" = 1.954678787305056; "H" = 1.574234825150171; "I" = 1.204555078994105; "J" = 0.08683365125180842; "M" = 0.7068601050726926 }
    16 = @{ "B" = 1.464709591687438; "C" = 0.5264786437560929; "D" = 0.5920757833877985; "E" = 0.2081799787770606; "G" = 1.931465175995157; "H" = 1.567671462601851; "I" = 1.200363175849986; "J" = 0.08697467182584262; "M" = 0.6862068216204307 }
    17 = @{ "B" = 1.422267706975958; "C" = 0.5113292619856225; "D" = 0.5897476290562054; "E" = 0.2078112492809296; "G" = 1.917577217561046; "H" = 1.563879172151928; "I" = 1.197991069557681; "J" = 0.08706975658021676; "M" = 0.673625382224543 }
    18 = @{ "B" = 1.397906142717659; "C" = 0.502636904556141; "D" = 0.5884363161121939; "E" = 0.2076082969625723; "G" = 1.909718277993449; "H" = 1.561783740708393; "I" = 1.196699729899841; "J" = 0.08712759475555032; "M" = 0.6664214062439697 }
    19 = @{ "B" = 1.389666301794534; "C" = 0.4996974469479483; "D" = 0.5879970960893104; "E" = 0.2075411481110194; "G" = 1.907079481274224; "H" = 1.561088975658834; "I" = 1.196275021586139; "J" = 0.08714771817784595; "M" = 0.6639878466305049 }
    20 = @{ "B" = 1.426780554695483; "C" = 0.5129397464990575; "D" = 0.5899925891867781; "E" = 0.2078495558982496; "G" = 1.919042244654008; "H" = 1.564273983160604; "I" = 1.19823601941863; "J" = 0.08705930878652524; "M" = 0.674961329443704 }
    21 = @{ "B" = 1.552075538161205; "C" = 0.5576857067910623; "D" = 0.5970293219894245; "E" = 0.2089949332779248; "G" = 1.960781677447301; "H" = 1.576001937218194; "I" = 1.205699248630395; "J" = 0.08679989521913001; "M" = 0.7122198393408183 }
    22 = @{ "B" = 1.634407814087069; "C" = 0.5871205769499852; "D" = 0.6018763459015872; "E" = 0.209825036202016; "G" = 1.989223183591349; "H" = 1.584436032670425; "I" = 1.211233193473788; "J" = 0.08665844311423498; "M" = 0.7368613416512204 }
    23 = @{ "B" = 1.590424800780454; "C" = 0.5713930601809238; "D" = 0.5992666593616036; "E" = 0.2093745235114923; "G" = 1.973936415646989; "H" = 1.579863961121447; "I" = 1.208219283197863; "J" = 0.08673136571746198; "M" = 0.7236830298483596 }
    24 = @{ "B" = 1.424740175265299; "C" = 0.5122115928311928; "D" = 0.5898817580790876; "E" = 0.2078322093361642; "G" = 1.91837951534086; "H" = 1.564095224986545; "I" = 1.198125052078694; "J" = 0.08706402235194588; "M" = 0.6743572565766272 }
    25 = @{ "B" = 1.24792897854303; "C" = 0.4491825153622813; "D" = 0.5808202410410672; "E" = 0.2065175407182984; "G" = 1.863378092139129; "H" = 1.550363009692518; "I" = 1.190025426017996; "J" = 0.0875431655972605; "M" = 0.6223944826080086 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
